$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply the per-cell text updates described by the diff.
# Every touched cell holds plain text (inlineStr) in the source workbook,
# including values that look numeric (e.g. "623.71", "70.187.72").
# A leading apostrophe forces Excel to store those as text instead of
# auto-converting them to numbers; ClearFormats() afterwards drops the
# resulting "quote prefix" cell style so formatting stays untouched.

$ws.Range('D2').Value = '70.187.72'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '3.738.44'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''623.71'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = '''180.94'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('D7').Value = '3.737.17'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '''0.535'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('E10').Value = '  +2.16%  '
$ws.Range('E11').Value = '  -5.02%  '
$ws.Range('E12').Value = '  -3.17%  '
$ws.Range('D13').Value = '''41.00'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.82%  '
$ws.Range('E14').Value = '  +1.81%  '
$ws.Range('D15').Value = '4.363.38'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('D16').Value = '3.741.63'
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('D17').Value = '70.164.83'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('D19').Value = '''7.60'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').Value = '''16.83'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').Value = '''505.84'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.51%  '
$ws.Range('D22').Value = '''9.37'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').Value = '''0.725'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.61%  '
$ws.Range('D24').Value = '''2.56'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('D25').Value = '''86.73'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.49%  '
$ws.Range('D26').Value = '''13.15'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.52%  '
$ws.Range('D27').Value = '''11.39'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.52%  '
$ws.Range('E28').Value = '  +23.49%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '''2.50'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('D32').Value = '''7.94'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.26%  '
$ws.Range('D33').Value = '''31.27'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.03%  '
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('E36').Value = '  +2.21%  '
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('D38').Value = '''0.137'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.11%  '
$ws.Range('E39').Value = '  -3.30%  '
$ws.Range('D40').Value = '''2.12'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -6.50%  '
$ws.Range('D41').Value = '''50.43'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('D42').Value = '''45.54'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').Value = '''427.40'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.15%  '
$ws.Range('D44').Value = '''8.75'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('D45').Value = '''2.88'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = '3.005.80'
$ws.Range('E46').Value = '  -4.49%  '
$ws.Range('E47').Value = '  -1.36%  '
$ws.Range('D48').Value = '''27.47'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.96%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '''2.52'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''136.79'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.46%  '
